$d = $word.ActiveDocument

# Hybrid bold + color highlighting for quantitative impact metrics.
# wdColor value for hex 2C3E50 (BGR-packed integer: R + G*256 + B*65536)
$metricColor = 5258796

function Highlight-Metrics {
    param(
        [int]$ParaIndex,
        [string[]]$Metrics
    )

    $para = $d.Paragraphs.Item($ParaIndex)
    $paraRange = $para.Range
    $searchStart = $paraRange.Start
    $searchEnd = $paraRange.End

    foreach ($metric in $Metrics) {
        $sub = $d.Range($searchStart, $searchEnd)
        $found = $sub.Find.Execute($metric, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $sub.Font.Bold = 1
            $sub.Font.Color = $metricColor
            $searchStart = $sub.End
        }
    }
}

# 1) "Discovered systematic race coding errors ... from 23% to 64%"
Highlight-Metrics 9 @("23%", "64%")

# 2) "Achieved 87% prediction accuracy ... from ±4.2% to ±2.1%"
Highlight-Metrics 11 @("87%", "71%", "±4.2%", "±2.1%")

# 3) "Wrote RFP and analyzed bids from 1,200 vendors ..."
Highlight-Metrics 31 @("1,200")

# 4) "... became the $400M Polling Consortium Database ... valued at $1B+"
Highlight-Metrics 46 @("$400M", "$1B")

# 5) "Algorithm reduced mapping costs by 73.5%, saving ... $4.7M"
Highlight-Metrics 63 @("73.5%", "$4.7M")

# 6) "Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%"
Highlight-Metrics 65 @("87%", "71%")

Write-Output "metrics highlighted"
